$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B:D range for rows 2-9 to 0 by default
$ws.Range("B2:D9").Value = 0

# Apply the two non-zero updated values from the diff
$ws.Range("C4").Value = -0.6449845920574532
$ws.Range("C8").Value = 0.7189829351727091
